# Add "2022-Q3" data to the 创业板/300566-激智科技 workbook:
#  1. Insert a new worksheet named "2022-Q3" right before the existing "2022-Q1"
#     sheet and fill it with the quarterly fund-holding table.
#  2. Insert a new row into the "总计" (summary) sheet for the 2022-Q3 figures,
#     pushing the existing quarters down and adding 2020-Q4 at the bottom.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet, positioned before "2022-Q1".
# ---------------------------------------------------------------------------
$refSheet = $wb.Worksheets.Item("2022-Q1")
$ws = $wb.Worksheets.Add($refSheet)
$ws.Name = "2022-Q3"

# Pull header + index-column formatting from an existing quarter sheet so the
# new tab looks like all the others (bold, bordered, centered style).
$fmtSource = $wb.Worksheets.Item("2021-Q4")
$fmtSource.Range("A1:H1").Copy($ws.Range("A1:H1"))
$fmtSource.Range("A2:A6").Copy($ws.Range("A2:A6"))

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Fund code / name / ratio columns are stored as plain text in this workbook
# (leading zeros in fund codes, fixed decimal strings) -- force text format
# before assigning so the engine doesn't "helpfully" coerce them to numbers.
$ws.Range("B2:G6").NumberFormat = "@"

$fundRows = @(
  @("160212", "国泰估值优势混合（LOF）A", "9.14", "94.29", "6.83", "0.6243", 5),
  @("162720", "广发创业板两年定期开放混合", "6.33", "94.20", "4.02", "0.2545", 8),
  @("007731", "民生加银持续成长混合A",     "3.22", "94.57", "5.70", "0.1835", 6),
  @("007732", "民生加银持续成长混合C",     "1.89", "94.57", "5.70", "0.1077", 6),
  @("016616", "国泰估值优势混合（LOF）C",  "0.00", "94.29", "6.83", "__NUM0__", 5)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
  $r = $i + 2
  $row = $fundRows[$i]
  $ws.Cells.Item($r, 1).Value = $i
  $ws.Cells.Item($r, 2).Value = $row[0]
  $ws.Cells.Item($r, 3).Value = $row[1]
  $ws.Cells.Item($r, 4).Value = $row[2]
  $ws.Cells.Item($r, 5).Value = $row[3]
  $ws.Cells.Item($r, 6).Value = $row[4]
  if ($row[5] -eq "__NUM0__") {
    $ws.Cells.Item($r, 7).NumberFormat = "General"
    $ws.Cells.Item($r, 7).Value = 0
  } else {
    $ws.Cells.Item($r, 7).Value = $row[5]
  }
  $ws.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: shift the quarter/count/value columns
#    down a row and add the new 2022-Q3 entry at the top (the row index in
#    column A is just a 0-based counter and is unaffected by the shift).
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

$summaryRows = @(
  @("2022-Q3", 5, 1.17),
  @("2022-Q1", 4, 0.98),
  @("2021-Q4", 12, 5.66),
  @("2021-Q3", 6, 0.63),
  @("2021-Q2", 6, 5.13),
  @("2021-Q1", 10, 8.1),
  @("2020-Q4", 6, 3.13)
)

for ($i = 0; $i -lt $summaryRows.Count; $i++) {
  $r = $i + 2
  $row = $summaryRows[$i]
  $totalWs.Cells.Item($r, 2).Value = $row[0]
  $totalWs.Cells.Item($r, 3).Value = $row[1]
  $totalWs.Cells.Item($r, 4).Value = $row[2]
}

# Row 8 (2020-Q4) is brand new -- copy the index-column style from the row
# above it before writing the new index value.
$totalWs.Range("A7").Copy($totalWs.Range("A8"))
$totalWs.Range("A8").Value = 6
